$d = $word.ActiveDocument

# Locate the "Part B" bullet that currently ends without a period.
$rng = $d.Range(0, $d.Content.End)
$found = $rng.Find.Execute( `
    "Take out the total number for the most 2 colors that are in the drawer", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target sentence"
}

# Collapse to the end of the found text (just before the _GoBack bookmark).
$rng.Collapse(0)
$periodStart = $rng.Start

# Insert the closing period as its own run (same rPr as the sentence before it).
$rng.InsertAfter(".")

# Briefly toggle Bold on just the new "." run so Word is forced to materialize
# it as a distinct <w:r> instead of folding it back into the previous run's text.
$periodRng = $d.Range($periodStart, $periodStart + 1)
$periodRng.Font.Bold = $true
$periodRng.Font.Bold = $false

# Split off a new list paragraph right after the period (inherits the same
# ListParagraph / numId 7, ilvl 1 formatting).
$afterPeriod = $d.Range($periodStart + 1, $periodStart + 1)
$afterPeriod.InsertParagraphAfter()

# Fill in the new paragraph's text.
$newParaTextPos = $periodStart + 2
$newParaRng = $d.Range($newParaTextPos, $newParaTextPos)
$newParaRng.InsertAfter( `
    "After taking out 16 socks, I will take out another 2, which for sure, I will have at least 1 matching pair of each color of socks.")

Write-Output "edit applied"
